$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 108. This shifts the existing rows 108-194 down to 109-195,
# preserving all of their data/formatting (including the date style on column D).
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with a new weekly price observation.
# All attributes mirror the record that was previously at row 108 (now row 109),
# except for the reporting date (D) and the reported volume (J).
$ws.Range("A108").Value = 10
$ws.Range("B108").Value = "Vega Modelo de Temuco"
$ws.Range("C108").Value = "La Araucanía"
$ws.Range("D108").Value = 44827
$ws.Range("E108").Value = 9
$ws.Range("F108").Value = 100114007
$ws.Range("G108").Value = "Jengibre"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 100
$ws.Range("K108").Value = 20000
$ws.Range("L108").Value = 20000
$ws.Range("M108").Value = 20000
$ws.Range("N108").Value = "$/caja 13 kilos"
$ws.Range("O108").Value = "Perú"
$ws.Range("P108").Value = 1538
$ws.Range("Q108").Value = 13
$ws.Range("R108").Value = "Hortaliza"
